$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44518
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 6000
$ws.Range("M2").Value = 6500
$ws.Range("O2").Value = 'Provincia de Diguillín'
$ws.Range("P2").Value = 260
$ws.Range("D3").Value = 44537
$ws.Range("K3").Value = 6500
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 6750
$ws.Range("O3").Value = 'Provincia de Diguillín'
$ws.Range("P3").Value = 270
$ws.Range("D4").Value = 44566
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 7500
$ws.Range("M4").Value = 7250
$ws.Range("O4").Value = 'Provincia de Diguillín'
$ws.Range("P4").Value = 290
$ws.Range("D5").Value = 44489
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 8500
$ws.Range("O5").Value = 'Región del Maule'
$ws.Range("P5").Value = 340
$ws.Range("D6").Value = 44167
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 8000
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 8500
$ws.Range("O6").Value = 'Región del Maule'
$ws.Range("P6").Value = 340
$ws.Range("D7").Value = 44482
$ws.Range("J7").Value = 120
$ws.Range("D8").Value = 44495
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 8000
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 8500
$ws.Range("O8").Value = 'Región del Maule'
$ws.Range("P8").Value = 340
$ws.Range("D9").Value = 44529
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 6000
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 6500
$ws.Range("O9").Value = 'Provincia de Diguillín'
$ws.Range("P9").Value = 260
$ws.Range("D10").Value = 44509
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 8500
$ws.Range("P10").Value = 340
$ws.Range("D11").Value = 44526
$ws.Range("D12").Value = 44540
$ws.Range("K12").Value = 6500
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 6750
$ws.Range("P12").Value = 270
$ws.Range("D13").Value = 44519
$ws.Range("J13").Value = 80
$ws.Range("D14").Value = 44536
$ws.Range("J14").Value = 80
$ws.Range("K14").Value = 6500
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 6750
$ws.Range("O14").Value = 'Provincia de Diguillín'
$ws.Range("P14").Value = 270
$ws.Range("D15").Value = 44166
$ws.Range("J15").Value = 56
$ws.Range("K15").Value = 7500
$ws.Range("L15").Value = 8000
$ws.Range("M15").Value = 7804
$ws.Range("O15").Value = 'Región de O''Higgins'
$ws.Range("P15").Value = 312
$ws.Range("D16").Value = 44476
$ws.Range("J16").Value = 160
$ws.Range("M16").Value = 7750
$ws.Range("O16").Value = 'Región del Maule'
$ws.Range("P16").Value = 310
$ws.Range("D17").Value = 44161
$ws.Range("J17").Value = 53
$ws.Range("K17").Value = 6500
$ws.Range("L17").Value = 7000
$ws.Range("M17").Value = 6764
$ws.Range("P17").Value = 271
$ws.Range("D18").Value = 44484
$ws.Range("J18").Value = 30
$ws.Range("K18").Value = 8500
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = 8750
$ws.Range("O18").Value = 'Región del Maule'
$ws.Range("P18").Value = 350
$ws.Range("D19").Value = 44488
$ws.Range("J19").Value = 60
$ws.Range("K19").Value = 8000
$ws.Range("L19").Value = 9000
$ws.Range("M19").Value = 8500
$ws.Range("O19").Value = 'Región del Maule'
$ws.Range("P19").Value = 340
$ws.Range("D21").Value = 44523
$ws.Range("K21").Value = 6000
$ws.Range("M21").Value = 6500
$ws.Range("P21").Value = 260
$ws.Range("D22").Value = 44491
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 8000
$ws.Range("L22").Value = 9000
$ws.Range("M22").Value = 8500
$ws.Range("O22").Value = 'Región del Maule'
$ws.Range("P22").Value = 340
$ws.Range("D23").Value = 44533
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = 6500
$ws.Range("L23").Value = 7000
$ws.Range("M23").Value = 6750
$ws.Range("O23").Value = 'Provincia de Diguillín'
$ws.Range("P23").Value = 270
$ws.Range("D24").Value = 44160
$ws.Range("K24").Value = 6500
$ws.Range("M24").Value = 6688
$ws.Range("O24").Value = 'Región de O''Higgins'
$ws.Range("P24").Value = 268
$ws.Range("D25").Value = 44159
$ws.Range("J25").Value = 42
$ws.Range("K25").Value = 6500
$ws.Range("L25").Value = 7000
$ws.Range("M25").Value = 6738
$ws.Range("P25").Value = 270
$ws.Range("D27").Value = 44511
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 7000
$ws.Range("L27").Value = 8000
$ws.Range("M27").Value = 7500
$ws.Range("O27").Value = 'Provincia de Diguillín'
$ws.Range("P27").Value = 300
$ws.Range("D28").Value = 44162
$ws.Range("J28").Value = 80
$ws.Range("K28").Value = 7000
$ws.Range("L28").Value = 8000
$ws.Range("M28").Value = 7562
$ws.Range("O28").Value = 'Región de O''Higgins'
$ws.Range("P28").Value = 302
$ws.Range("D29").Value = 44504
$ws.Range("J29").Value = 60
$ws.Range("K29").Value = 8000
$ws.Range("L29").Value = 9000
$ws.Range("M29").Value = 8500
$ws.Range("P29").Value = 340
$ws.Range("D30").Value = 44512
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 7000
$ws.Range("L30").Value = 8000
$ws.Range("M30").Value = 7500
$ws.Range("P30").Value = 300
$ws.Range("D31").Value = 44522
$ws.Range("D32").Value = 44524
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 6000
$ws.Range("L32").Value = 7000
$ws.Range("M32").Value = 6500
$ws.Range("O32").Value = 'Provincia de Diguillín'
$ws.Range("P32").Value = 260
$ws.Range("D33").Value = 44515
$ws.Range("K33").Value = 7000
$ws.Range("L33").Value = 8000
$ws.Range("M33").Value = 7500
$ws.Range("P33").Value = 300
$ws.Range("D34").Value = 44165
$ws.Range("J34").Value = 38
$ws.Range("L34").Value = 8500
$ws.Range("M34").Value = 8263
$ws.Range("P34").Value = 331
$ws.Range("D35").Value = 44553
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 6500
$ws.Range("L35").Value = 7000
$ws.Range("M35").Value = 6750
$ws.Range("O35").Value = 'Provincia de Diguillín'
$ws.Range("P35").Value = 270
$ws.Range("D37").Value = 44517
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = 6000
$ws.Range("M37").Value = 6500
$ws.Range("P37").Value = 260
$ws.Range("D38").Value = 44487
$ws.Range("J38").Value = 30
$ws.Range("L38").Value = 8000
$ws.Range("M38").Value = 8000
$ws.Range("P38").Value = 320
$ws.Range("D39").Value = 44487
$ws.Range("I39").Value = 'Segunda'
$ws.Range("J39").Value = 30
$ws.Range("K39").Value = 9000
$ws.Range("L39").Value = 9000
$ws.Range("M39").Value = 9000
$ws.Range("O39").Value = 'Región del Maule'
$ws.Range("P39").Value = 360
$ws.Range("D40").Value = 44466
$ws.Range("J40").Value = 60
$ws.Range("K40").Value = 11000
$ws.Range("L40").Value = 12000
$ws.Range("M40").Value = 11500
$ws.Range("O40").Value = 'Región de O''Higgins'
$ws.Range("P40").Value = 340
"Applied " + 206 + " cell updates"
